$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows for Sheet1 (graph size benchmarks) -- values changed, and last two
# rows (previous 50001 entries) removed, shrinking the table from A1:E13 to A1:E11.
$data = @(
    @(20,    2, "0s",      0.5, 2),
    @(20,    3, "0s",      0.5, 3),
    @(200,   2, "0,01s",   0.5, 2),
    @(200,   3, "0,01s",   0.5, 3),
    @(2000,  2, "0,26s",   0.5, 2),
    @(2000,  3, "0,2s",    0.5, 3),
    @(20000, 2, "30,43s",  0.5, 2),
    @(20000, 3, "31,2s",   0.5, 3),
    @(50000, 2, "291,87s", 0.5, 2),
    @(50000, 3, "286,46s", 0.5, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Remove the old two trailing rows (12 and 13), which no longer exist in the new table.
$ws.Range("A12:E13").EntireRow.Delete()
